# "add tabel format baru"
# Bumps the table numbering (Tabel 4.2.3/4.2.4./4.2.5. -> 4.2.5/4.2.6./4.2.7.)
# and the reporting year (2020 -> 2021) for the Kecamatan Iwoimendaa sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table number headers (row 1) ---
$ws.Range("H1").Value = "Tabel 4.2.5"
$ws.Range("P1").Value = "Tabel 4.2.6."
$ws.Range("W1").Value = "Tabel 4.2.7."

# --- Table titles: bump year 2020 -> 2021 (row 1, Indonesian captions) ---
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Iwoimendaa. 2021"
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Iwoimendaa, 2021"
$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Iwoimendaa, 2021"
$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Iwoimendaa, 2021"

# --- Table titles: bump year 2020 -> 2021 (row 2, English captions) ---
$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in Iwoimendaa Subdistrict, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village Iwoimendaa Subdistrict, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in Iwoimendaa Subdistrict, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in Iwoimendaa Subdistrict, 2021"

# --- Reset the saved view back to the sheet's top-left / default selection ---
$ws.Range("A1").Select()
